# edit.ps1 — apply the "Project Management Plan.docx" revision described
# by the commit: grammar/spelling fixes, figure field fix, image tweak,
# bookmark relocation, and section margin changes.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARN: not found: $find"
    }
}

# 1. "configured that user." -> "configured by that user."
Replace-Text "configured that user." "configured by that user."

# 2. "will allows users" -> "will allow users"
Replace-Text "The GUI of the application will allows users" "The GUI of the application will allow users"

# 3. "Calibri Body, 12, Bold" (Subsection) -> "Calibri Body, 11, Bold"
Replace-Text "Subsection: Calibri Body, 12, Bold" "Subsection: Calibri Body, 11, Bold"

# 4. "each major sections should" -> "each major section should"
Replace-Text "Format: The heading of each major sections should" "Format: The heading of each major section should"

# 5. "google hangouts" -> "Google Hangouts"
Replace-Text "through google hangouts if" "through Google Hangouts if"

# 6. "due to ability of" -> "due to the ability of"
Replace-Text "Agile is very flexible due to ability of customers" "Agile is very flexible due to the ability of customers"

# 7. "problems you had, what you plan" -> "problems you had, and what you plan"
Replace-Text "problems you had, what you plan to work on" "problems you had, and what you plan to work on"

# 8. "using Visual Studio using" -> "using Microsoft Visual Studio using"
Replace-Text "The JARVIS Emulator will be programmed using Visual Studio using the .NET framework" "The JARVIS Emulator will be programmed using Microsoft Visual Studio using the .NET framework"

# 9. "will be use in the process" -> "will be used in the process"
Replace-Text "(WPF) will be use in the process" "(WPF) will be used in the process"

# 10. "Windows Operating System" (end of paragraph, before line break) gets a period
Replace-Text "This application will work on the Windows Operating System" "This application will work on the Windows Operating System."

# 11. "we be using TFS" -> "we will be using TFS"
Replace-Text "r (TFS), we be using TFS" "r (TFS), we will be using TFS"

# 12. "guarantee the stability of the repository." -> "guarantee the repository's stability."
Replace-Text "guarantee the stability of the repository." "guarantee the repository’s stability."

# 13. "He must be sure to communicate" -> "He must communicate"
Replace-Text "He must be sure to communicate with the rest of the team" "He must communicate with the rest of the team"

# 14. "these type of situations" -> "these types of situations"
Replace-Text "handling these type of situations." "handling these types of situations."

# 15. add trailing "." after "...calendar events) "
Replace-Text "bother some users (maybe remembering birthdays or calendar events) " "bother some users (maybe remembering birthdays or calendar events) ."

# 16. Security paragraph rewording
Replace-Text " Security will be important, since certain information that the user may input on their profile might be of interest such as phishing and targeted ads, redirected websites that are not the predetermined ones." "Security will be important, since certain information that the user may input on his or her profile might be of interest to phishing, targeted ads and redirected websites that are not the predetermined ones."

# 17. "Develop website APIs modules" -> "Develop website API modules"
Replace-Text "Develop website APIs modules" "Develop website API modules"

# 18. Section page-margin tweak: top 63pt->54pt (1260->1080 twips), bottom 54pt->49.5pt (1080->990 twips)
$ps = $d.Sections(1).PageSetup
$ps.TopMargin = 54
$ps.BottomMargin = 49.5

Write-Host "done part 1"
